# Update numeric values (column F = sales/count, column G = price) across
# the four worksheets of the 广州-漫展信息 workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 228
$ws1.Range("F5").Value = 1627
$ws1.Range("F7").Value = 623
$ws1.Range("F8").Value = 133
$ws1.Range("F9").Value = 598
$ws1.Range("F10").Value = 54
$ws1.Range("F11").Value = 101

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("G4").Value = 149
$ws2.Range("F11").Value = 14

# --- Sheet "本地生活" (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 1941
$ws3.Range("F5").Value = 43

# --- Sheet "全部类型" (All types, combined view) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1941
$ws4.Range("F6").Value = 43
$ws4.Range("G10").Value = 149
$ws4.Range("F12").Value = 228
$ws4.Range("F16").Value = 1627
$ws4.Range("F20").Value = 14
$ws4.Range("F21").Value = 623
$ws4.Range("F23").Value = 133
$ws4.Range("F24").Value = 598
$ws4.Range("F25").Value = 54
$ws4.Range("F27").Value = 101
